$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.923.54'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.038.03'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.97'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.46'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.375'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0769'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.36'
$ws.Range("E12").Value = '  -2.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.878'
$ws.Range("E13").Value = '  +8.69%  '
$ws.Range("D14").Value = '2.335.03'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.62'
$ws.Range("E15").Value = '  +2.12%  '
$ws.Range("D16").Value = '2.037.87'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("E17").Value = '  +8.88%  '
$ws.Range("D18").Value = '36.898.05'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.52'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.37'
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.29'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("E24").Value = '  +3.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.81'
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("E27").Value = '  -3.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.89'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.43'
$ws.Range("E29").Value = '  +16.07%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.70'
$ws.Range("E32").Value = '  +6.67%  '
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0864'
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("E36").Value = '  +6.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.24'
$ws.Range("E37").Value = '  +1.21%  '
$ws.Range("E38").Value = '  -3.05%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("E40").Value = '  +1.18%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0964'
$ws.Range("E42").Value = '  -11.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.13'
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.19'
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.95'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("D46").Value = '1.293.02'
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.76'
$ws.Range("E48").Value = '  +10.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.70'
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").Value = '2.222.44'
$ws.Range("E51").Value = '  -0.63%  '
